# "filter error in worksheet"
#
# The commit adds a new "error" column (N) alongside the thickness values
# already present in column M (rows 39-59), with a header label "error" in
# N38 (mirroring the existing PARAMETER/VALUE/ERROR BAR style layout used
# elsewhere on the sheet), repositions/resizes the chart slightly, turns on
# explicit page-setup (paper size + orientation) and updates the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "error" column (N) data, rows 38-59
# ---------------------------------------------------------------------
$ws.Range("N38").Value = "error"

$ws.Range("N39").Value = 0.018
$ws.Range("N40").Value = 0.022
$ws.Range("N41").Value = 0.018
$ws.Range("N42").Value = 0.021
$ws.Range("N43").Value = 0.022
$ws.Range("N44").Value = 0.032
$ws.Range("N45").Value = 0.026
$ws.Range("N46").Value = 0.022
$ws.Range("N47").Value = 0.039
$ws.Range("N48").Value = 0.035
$ws.Range("N49").Value = 0.029
$ws.Range("N50").Value = 0.032
$ws.Range("N51").Value = 0.047
$ws.Range("N52").Value = 0.039
$ws.Range("N53").Value = 0.051
$ws.Range("N54").Value = 0.049
$ws.Range("N55").Value = 0.04
$ws.Range("N56").Value = 0.055
$ws.Range("N57").Value = 0.044
$ws.Range("N58").Value = 0.057
$ws.Range("N59").Value = 0.061

# ---------------------------------------------------------------------
# 2. Reposition / resize the scatter chart (matches the twoCellAnchor
#    from/to change: from N37+24pt,3.3pt to V51,0pt+9.3pt)
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 894.6728515625
$co.Top = 543.3
$co.Width = 443.5
$co.Height = 216

# ---------------------------------------------------------------------
# 3. Page setup: paper size 9 (A4) and portrait orientation
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
$ps.HorizontalDpi = 300
$ps.VerticalDpi = 300

# ---------------------------------------------------------------------
# 4. Update the active selection / view
# ---------------------------------------------------------------------
$ws.Range("O56").Select()
